$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking / date-looking strings are preserved verbatim
$ws.Range("C2:O4").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "226.6"
$ws.Range("D2").Value = "233.9"
$ws.Range("E2").Value = "235.2"
$ws.Range("F2").Value = "229.4"
$ws.Range("G2").Value = "232.2"
$ws.Range("H2").Value = "PRECISION INSTRUMENTS"
$ws.Range("I2").Value = "145119"
$ws.Range("J2").Value = "C3FR250F"
$ws.Range("K2").Value = "2025-03-06"
$ws.Range("L2").Value = "2026-03-06"
$ws.Range("M2").Value = "03"
$ws.Range("N2").Value = "KALTIRE 089"
$ws.Range("O2").Value = "780-228-5970"

# Row 3
$ws.Range("C3").Value = "147.5"
$ws.Range("D3").Value = "144.3"
$ws.Range("E3").Value = "146.5"
$ws.Range("F3").Value = "145.7"
$ws.Range("G3").Value = "145.7"
$ws.Range("H3").Value = "PRECISION INSTRUMENTS"
$ws.Range("I3").Value = "145119"
$ws.Range("J3").Value = "C3FR250F"
$ws.Range("K3").Value = "2025-03-06"
$ws.Range("L3").Value = "2026-03-06"
$ws.Range("M3").Value = "03"
$ws.Range("N3").Value = "KALTIRE 089"
$ws.Range("O3").Value = "780-228-5970"

# Row 4
$ws.Range("C4").Value = "72.1"
$ws.Range("D4").Value = "70.2"
$ws.Range("E4").Value = "72.1"
$ws.Range("F4").Value = "71.7"
$ws.Range("G4").Value = "71.0"
$ws.Range("H4").Value = "PRECISION INSTRUMENTS"
$ws.Range("I4").Value = "145119"
$ws.Range("J4").Value = "C3FR250F"
$ws.Range("K4").Value = "2025-03-06"
$ws.Range("L4").Value = "2026-03-06"
$ws.Range("M4").Value = "03"
$ws.Range("N4").Value = "KALTIRE 089"
$ws.Range("O4").Value = "780-228-5970"
